$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 3-5 entirely (table shrinks from 5 data rows to 2)
$ws.Rows("3:5").Delete()

# Fill in column A (office names) first so the shared-string table
# gets them allocated before the region names.
$ws.Range("A1").Value = "office0"
$ws.Range("A2").Value = "office1"

# Region columns (C and D)
$ws.Range("C1").Value = "region0"
$ws.Range("D1").Value = "region1"
$ws.Range("C2").Value = "region2"
$ws.Range("D2").Value = "region1"

# Demand values in column B, right-aligned
$ws.Range("B1").Value = 50
$ws.Range("B2").Value = 50
$ws.Range("B1:B2").HorizontalAlignment = -4152

# The text cells lose their old centered style and go back to the
# default (unstyled) cell format
$ws.Range("A1:A2").HorizontalAlignment = 1
$ws.Range("C1:D2").HorizontalAlignment = 1

# View settings
$excel.ActiveWindow.Zoom = 115
[void]$ws.Range("D5").Select()
